{"js": "// Revert the earlier change that added a comment (\"similar to the first\n// point\") anchored to the text:\n//   \"Cybersecurity measures cannot be copy-pasted from larger business to\n//    smaller ones. (ibid. 6)\"\n// This removes the comment (and its anchor marks: commentRangeStart,\n// commentRangeEnd, commentReference) from the document entirely.\n\nconst comments = context.document.body.getComments();\ncomments.load(\"items\");\nawait context.sync();\n\nfor (let i = comments.items.length - 1; i >= 0; i--) {\n  comments.items[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# Revert the earlier change that added a comment (\"similar to the first\n# point\") anchored to the text:\n#   \"Cybersecurity measures cannot be copy-pasted from larger business to\n#    smaller ones. (ibid. 6)\"\n# This removes the comment (and its anchor marks: commentRangeStart,\n# commentRangeEnd, commentReference) from the document entirely.\n\n$d = $word.ActiveDocument\n\nfor ($i = $d.Comments.Count; $i -ge 1; $i--) {\n    $d.Comments($i).Delete()\n}\n"}
